$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "At-Risk to PB" row (row 7) is removed entirely; Excel shifts rows
# 8:28 up to 7:27 (formulas/values/styles move with the row, dimension
# shrinks from A1:S28 to A1:S27, and the now-orphaned "At-Risk to PB"
# shared string is dropped from sharedStrings.xml on save).
$ws.Rows.Item(7).Delete()

# The row that was "At-risk" (now still row 5 post-delete) is relabeled.
$ws.Range("A5").Value = "Negative Equity"

# Reflect the author's resulting view/selection state: scrolled down with
# A13 as the active cell.
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A13").Select()
